$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns keep their original text formatting
# so numeric-looking values (e.g. "1.01") are not auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.519.89'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.628.15'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.93'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.99%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.01'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.854.36'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.629.67'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.45%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.49'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.546.21'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.59'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.34'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.55%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.23'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.01%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.53'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.217.81'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0173'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.31%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.25%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.68%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.765.02'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.42'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.57'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.90'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0512'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.65'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.60%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.06%  '
